# Generate Report for Handoff
#
# The three tracked files (identified by their e2e markdown file name) are
# re-sorted alphabetically across the "Overview", "zh-cn" and "de-de"
# sheets, and the row belonging to
# "b7072dda-8467-44f1-81dd-51f4233aa325.md" is refreshed with a new
# handoff status/timestamp (it moves from "Handed back: in sync with
# en-US" to "Ready for handoff").
#
# Helper: write a value into a cell AND keep any hyperlink anchored to
# that cell in sync (the model does not auto-sync hyperlink display text
# when a cell's value is overwritten).
function Set-CellAndHyperlink {
    param(
        $Worksheet,
        [string]$CellRef,
        [string]$NewValue
    )

    $Worksheet.Range($CellRef).Value = $NewValue

    foreach ($h in $Worksheet.Hyperlinks) {
        if ($h.Range.Address() -eq ('$' + ($CellRef -replace '(\d+)', '$$$1'))) {
            $h.TextToDisplay = $NewValue
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $ws "A2" "ffffc07846bf-1d80-4ff4-983e-f1882168d505.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-19 17:05:25"

Set-CellAndHyperlink $ws "A3" "ffffffc683b9b7-85e0-4ba8-bdbc-c71846d81638.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-19 17:05:25"

Set-CellAndHyperlink $ws "A4" "b7072dda-8467-44f1-81dd-51f4233aa325.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-19 17:10:28"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $ws "A2" "ffffc07846bf-1d80-4ff4-983e-f1882168d505.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $ws "D2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-19 17:05:18"
Set-CellAndHyperlink $ws "F2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.md"
Set-CellAndHyperlink $ws "G2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-19 17:05:57"
$ws.Range("J2").Value = "Include"

Set-CellAndHyperlink $ws "A3" "ffffffc683b9b7-85e0-4ba8-bdbc-c71846d81638.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $ws "D3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-19 17:05:18"
Set-CellAndHyperlink $ws "F3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.md"
Set-CellAndHyperlink $ws "G3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-19 17:05:57"
$ws.Range("J3").Value = "Include"

Set-CellAndHyperlink $ws "A4" "b7072dda-8467-44f1-81dd-51f4233aa325.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
Set-CellAndHyperlink $ws "D4" "b7072dda-8467-44f1-81dd-51f4233aa325.3f8e7d29752cc1767baacbe738c0c70a2d4f1246.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-19 17:10:20"
Set-CellAndHyperlink $ws "F4" "b7072dda-8467-44f1-81dd-51f4233aa325.md"
Set-CellAndHyperlink $ws "G4" "b7072dda-8467-44f1-81dd-51f4233aa325.3f8e7d29752cc1767baacbe738c0c70a2d4f1246.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-19 17:09:39"
$ws.Range("J4").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $ws "A2" "ffffc07846bf-1d80-4ff4-983e-f1882168d505.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $ws "D2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.de-de.xlf"
$ws.Range("E2").Value = "2016-03-19 17:05:25"
Set-CellAndHyperlink $ws "F2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.md"
Set-CellAndHyperlink $ws "G2" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.de-de.xlf"
$ws.Range("H2").Value = "2016-03-19 17:06:14"
$ws.Range("J2").Value = "Include"

Set-CellAndHyperlink $ws "A3" "ffffffc683b9b7-85e0-4ba8-bdbc-c71846d81638.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $ws "D3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.de-de.xlf"
$ws.Range("E3").Value = "2016-03-19 17:05:25"
Set-CellAndHyperlink $ws "F3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.md"
Set-CellAndHyperlink $ws "G3" "6ee46d3b-90f0-4252-81ec-5208d71c7155.8ae29b0051948772396be05d69acdb0d9171e801.de-de.xlf"
$ws.Range("H3").Value = "2016-03-19 17:06:14"
$ws.Range("J3").Value = "Include"

Set-CellAndHyperlink $ws "A4" "b7072dda-8467-44f1-81dd-51f4233aa325.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
Set-CellAndHyperlink $ws "D4" "b7072dda-8467-44f1-81dd-51f4233aa325.3f8e7d29752cc1767baacbe738c0c70a2d4f1246.de-de.xlf"
$ws.Range("E4").Value = "2016-03-19 17:10:28"
Set-CellAndHyperlink $ws "F4" "b7072dda-8467-44f1-81dd-51f4233aa325.md"
Set-CellAndHyperlink $ws "G4" "b7072dda-8467-44f1-81dd-51f4233aa325.3f8e7d29752cc1767baacbe738c0c70a2d4f1246.de-de.xlf"
$ws.Range("H4").Value = "2016-03-19 17:09:53"
$ws.Range("J4").Value = "Include"
